# Mater update 22 nov 2020
# Re-sort the "Item Name" / "UOM" pairs for a few brand groups (Dinafex,
# Etorix, Ketonic, Kynol, Zithrox) so items within each brand are in
# alphabetical/ascending order, keeping Item Name and UOM values in sync.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Dinafex group (rows 3-5): reorder to 60mg, 120mg, 180mg
$ws.Range("D3").Value = "Dinafex 60mg Tablet"
$ws.Range("D5").Value = "Dinafex 180mg Tablet"

# Etorix group (rows 8-9): swap 90mg Tablet and 60mg Tablet - 40's
$ws.Range("D8").Value = "Etorix 60mg Tablet - 40's"
$ws.Range("E8").Value = "40's"
$ws.Range("D9").Value = "Etorix 90mg Tablet"
$ws.Range("E9").Value = "30's"

# Ketonic group (rows 14-16): reorder to Injection - 4's, 10mg Tablet, 30mg Injection
$ws.Range("D14").Value = "Ketonic 30mg IM/IV Injection - 4's"
$ws.Range("E14").Value = "4's"
$ws.Range("D15").Value = "Ketonic 10mg Tablet"
$ws.Range("E15").Value = "20's"
$ws.Range("D16").Value = "Ketonic 30mg Injection"
$ws.Range("E16").Value = "5 's"

# Kynol group (rows 17-18): swap 100mg Capsule and 200mg Capsule
$ws.Range("D17").Value = "Kynol TR 200mg Capsule"
$ws.Range("E17").Value = "30 's"
$ws.Range("D18").Value = "Kynol TR 100mg Capsule"
$ws.Range("E18").Value = "50 's"

# Zithrox group (rows 24,26): swap 15ml Suspension and 30ml Dry Suspension
$ws.Range("D24").Value = "Zithrox 30ml Dry Suspension"
$ws.Range("E24").Value = "30ml"
$ws.Range("D26").Value = "Zithrox 15ml Suspension"
$ws.Range("E26").Value = "15 ml"
